# newYearProject_api_0.8v.xlsx
# refactor: Letter, User의 create_at 추가
#
# Adds a new "get single letter by id" endpoint:
#   - sheet "request"  (sheet1): new row documenting GET api/letter/{uuid}/getLetter/{id}
#   - sheet "response" (sheet2): existing 친구페이지 sample response gains a "wish" field,
#       the paging endpoint description is clarified to "전체 조회", and a new row documents
#       the response body for the single-letter lookup.

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("request")
$ws2 = $wb.Worksheets.Item("response")

# ---------------------------------------------------------------------------
# Sheet "response" text edits first
# ---------------------------------------------------------------------------

$letterListJson = @"
{
"resultCode": "SUCCESS",
"result":{
{
"nickName":"sujan",
"wish":"다이어트",
"money":10000
}
}
"@

$ws2.Cells.Item(11, 3).Value = $letterListJson
$ws2.Rows.Item(11).RowHeight = 153

$letterDetailJson = @"
{
"resultCode": "SUCCESS",
	"result":{
		"author": "수진이",
		"content": "안녕? 새해복마니바다",
		"money": 10000
	}
}
"@

$ws2.Cells.Item(18, 3).Value = $letterDetailJson
$ws2.Cells.Item(18, 4).Value = "편지 상세 조회"

$ws2.Cells.Item(17, 4).Value = "편지 전체 조회, 페이징적용"

# ---------------------------------------------------------------------------
# Sheet "request": new row 10 -> single letter detail lookup
# ---------------------------------------------------------------------------

$ws1.Cells.Item(10, 3).Value = "api/letter/{uuid}/getLetter/{id}"
$ws1.Cells.Item(10, 5).Value = "편지 상세 조회"
$ws1.Cells.Item(10, 1).Value = "request"
$ws1.Cells.Item(10, 2).Value = "GET"

# yellow highlight on the URL column, matching the rest of the table
$ws1.Cells.Item(10, 3).Interior.Color = 65535
# wrap the long description text, matching the rest of the table
$ws1.Cells.Item(10, 5).WrapText = $true

# ---------------------------------------------------------------------------
# Sheet "request" row 9 keeps the same text, just shifts shared-string index
# (no-op in value, kept here for clarity/robustness)
# ---------------------------------------------------------------------------

$ws1.Cells.Item(9, 3).Value = "api/letter/{uuid}/getLetter"
$ws1.Cells.Item(9, 5).Value = "편지 목록 조회"

# ---------------------------------------------------------------------------
# Sheet "response" row 18: finish filling in A/B columns
# ---------------------------------------------------------------------------

$ws2.Cells.Item(18, 1).Value = "response"
$ws2.Cells.Item(18, 2).Value = 200

$ws2.Cells.Item(18, 3).WrapText = $true
$ws2.Rows.Item(18).RowHeight = 136

# ---------------------------------------------------------------------------
# Selections: sheet2 cell C18 selected, but sheet1 (request) stays the active
# tab with B10 selected -- so select sheet1 *last*.
# ---------------------------------------------------------------------------

$ws2.Range("C18").Select()
$ws1.Range("B10").Select()
